# FFM FUNCIONAL CON TENSORFLOW
# - Rename "ORTABLA" sheet to "XORTABLA"
# - Update its truth table so row 5 (1,1) -> C5 = 0 (XOR instead of OR)
# - Make XORTABLA the active sheet/tab with selection on E5

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ORTABLA")
$ws.Name = "XORTABLA"
$ws.Range("C5").Value = 0

$ws.Activate()
$ws.Range("E5").Select()
